$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    3  = -3
    6  = 1
    10 = 3
    13 = -2
    14 = 6
    19 = 0
    20 = -2
    26 = -3
    29 = 1
    30 = 1
    31 = -3
    33 = -1
    34 = -2
    40 = -7
    41 = -2
    47 = -1
    48 = -2
    50 = 2
    52 = -9
    54 = 0
    56 = -8
    57 = -3
    60 = -5
    61 = -6
    62 = -5
    64 = -2
    65 = -4
    67 = 4
    68 = 5
    69 = 6
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
